$p = $ppt.ActivePresentation

# --- Slide 13: "Rectangle 6" diagram box ---
$s13 = $p.Slides.Item(13)
$shp13 = $s13.Shapes.Item(4)

# Reposition / resize (height cy stays the same: 3779817 EMU)
$shp13.Left = 102.0                  # 1295400 EMU
$shp13.Top = 64.48574803149606       # 818969 EMU
$shp13.Width = 522.0                 # 6629400 EMU

$tr13 = $shp13.TextFrame.TextRange

$para4 = $tr13.Paragraphs(4,1)
$para4.Characters(1, $para4.Length).Text = "  Destination UDP Port            /  \      Network Programming Label"

$para5 = $tr13.Paragraphs(5,1)
$para5.Characters(1, $para5.Length).Text = "  Measurement Protocol           /    \     Timestamp2 Offset"

$para6 = $tr13.Paragraphs(6,1)
$para6.Characters(1, $para6.Length).Text = "  PLM Type                      /      \    Timestamp Format"

# The shape auto-fits its text box (spAutoFit); editing the text above
# recalculates Height, so pin it back to the unchanged target value.
$shp13.Height = 297.6234   # 3779817 EMU (unchanged by the edit)

# --- Slide 5: "Rectangle 2" diagram box ---
$s5 = $p.Slides.Item(5)
$shp5 = $s5.Shapes.Item(4)
$tr5 = $shp5.TextFrame.TextRange
$para5_3 = $tr5.Paragraphs(3,1)
$para5_3.Characters(1, $para5_3.Length).Text = " |   R1  |====================||  R5   |"

# --- Slide 8: "Rectangle 2" diagram box ---
$s8 = $p.Slides.Item(8)
$shp8 = $s8.Shapes.Item(4)
$tr8 = $shp8.TextFrame.TextRange
$para8_3 = $tr8.Paragraphs(3,1)
$para8_3.Characters(1, $para8_3.Length).Text = "   |   R1  |====================||  R5   |"
